$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.332.92"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.68%  "
$ws.Range("D3").Value = "'2.656.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.50%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'607.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("D6").Value = "'152.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.92%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("E8").Value = "  +1.09%  "
$ws.Range("E9").Value = "  +1.66%  "
$ws.Range("E10").Value = "  +6.60%  "
$ws.Range("E11").Value = "  -0.59%  "
$ws.Range("E12").Value = "  -0.67%  "
$ws.Range("D13").Value = "'28.22"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.36%  "
$ws.Range("D14").Value = "'3.135.64"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.57%  "
$ws.Range("D15").Value = "'64.298.34"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.92%  "
$ws.Range("D16").Value = "'0.0000148"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.76%  "
$ws.Range("D17").Value = "'2.659.69"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.68%  "
$ws.Range("D18").Value = "'12.15"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.87%  "
$ws.Range("D19").Value = "'4.64"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.65%  "
$ws.Range("D20").Value = "'349.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.01%  "
$ws.Range("D21").Value = "'6.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.66%  "
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").Value = "'5.58"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("D24").Value = "'66.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.34%  "
$ws.Range("E25").Value = "  +13.32%  "
$ws.Range("D26").Value = "'9.40"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +8.66%  "
$ws.Range("E27").Value = "  +4.38%  "
$ws.Range("D28").Value = "'8.19"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.93%  "
$ws.Range("D29").Value = "'553.19"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.57%  "
$ws.Range("D30").Value = "'0.165"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.49%  "
$ws.Range("E31").Value = "  -0.09%  "
$ws.Range("E32").Value = "  +1.47%  "
$ws.Range("D33").Value = "'0.0₃0864"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.88%  "
$ws.Range("E34").Value = "  -1.13%  "
$ws.Range("E35").Value = "  +6.14%  "
$ws.Range("D36").Value = "'168.72"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.96%  "
$ws.Range("E37").Value = "  +0.87%  "
$ws.Range("B38").Value = "FirstDigitalUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D38").Value = "'1.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "'1.97"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.15%  "
$ws.Range("D40").Value = "'19.43"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.82%  "
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("D42").Value = "'166.57"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.07%  "
$ws.Range("D43").Value = "'40.17"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.17%  "
$ws.Range("E44").Value = "  +2.86%  "
$ws.Range("D45").Value = "'0.0579"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.43%  "
$ws.Range("D46").Value = "'22.01"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.80%  "
$ws.Range("E47").Value = "  +0.30%  "
$ws.Range("E48").Value = "  +15.35%  "
$ws.Range("E49").Value = "  +3.39%  "
$ws.Range("E50").Value = "  +0.75%  "
$ws.Range("D51").Value = "'19.09"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.80%  "
